# clustTools and caretExtra compatibility
# Remove the "Brier skill score" column (column F) from the two supplementary
# tables that report PTS cluster sensitivity/specificity, since that metric
# is no longer produced by the updated clustTools/caretExtra packages.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Supplementary Table S11", "Supplementary Table S12")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(6).Delete()
}
